# feature: 产品页面支持 IPV6 主机导入 (closed #776)
#
# Adds a new "Addressing mode(Optional)" column (column O) to the node-
# import template sheet, with "Dynamic" / "Static" values on the sample
# rows, and reworks the "Speed limit M/s(Optional)" header cell so the
# zero-width-space glyphs in the middle are tagged with an explicit
# (Japanese) font run, matching upstream's rich-text split.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column O: "Addressing mode(Optional)" with Dynamic/Static values
# ---------------------------------------------------------------------
$ws.Range("O1").Value = "Addressing mode(Optional)"
# Write row 3 ("Static") before row 2 ("Dynamic") so the shared-string
# table picks up the same index order as upstream (33=Static, 34=Dynamic).
$ws.Range("O3").Value = "Static"
$ws.Range("O2").Value = "Dynamic"

# Match the column width used upstream (closest value the engine's
# character-width quantization can reach to 27.125).
$ws.Columns("O").ColumnWidth = 26.4

# ---------------------------------------------------------------------
# 2. Split the "Speed <zwsp><zwsp>limit M/s(Optional)" header (N1) into
#    a dedicated run for the zero-width-space pair using an explicit
#    Japanese font, and rejoin "limit M/s" with "(Optional)".
# ---------------------------------------------------------------------
$zwsp = [char]0x200b
$speedCell = $ws.Range("N1")
$speedCell.Value = "Speed " + $zwsp + $zwsp + "limit M/s(Optional)"
$speedCell.Characters(7, 2).Font.Name = "MS Gothic"
$speedCell.Characters(7, 2).Font.Size = 12

# ---------------------------------------------------------------------
# 3. Update the active selection to match the saved view (G17, no
#    frozen/scrolled topLeftCell override).
# ---------------------------------------------------------------------
$ws.Range("G17").Select()
